$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new ingredient row: hot dog, id 9, quantity 0
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "hot dog "
$ws.Range("C10").Value = 0

# Widen column B to fit the ingredient names (closest achievable to 22.85546875 chars)
$ws.Columns.Item(2).ColumnWidth = 22

# Update selection to match the new active cell
$ws.Range("C10").Select()
